# Update column G (K) values to reflect newly computed strikeout counts,
# replacing the old 'Strike#' derived values with the regenerated 'K' values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 3
    3 = 0
    4 = 1
    5 = 0
    6 = 0
    7 = 1
    8 = 0
    9 = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    19 = 2
    20 = 0
    21 = 3
    22 = 1
    23 = 3
    24 = 0
    25 = 1
    26 = 2
    27 = 0
    28 = 2
    29 = 1
    30 = 0
    31 = 0
    32 = 2
    33 = 1
    34 = 0
    35 = 3
    36 = 1
    37 = 1
    38 = 2
    39 = 2
    40 = 1
    41 = 0
    42 = 1
    43 = 0
    44 = 2
    45 = 0
    46 = 0
    47 = 2
    48 = 0
    49 = 0
    50 = 0
    51 = 2
    52 = 1
    53 = 0
    54 = 2
    55 = 0
    56 = 1
    57 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
